$wb = $excel.ActiveWorkbook

# --- Sheet 1: studyData ---
$ws1 = $wb.Worksheets.Item("studyData")

# disorder (column X) for study row 2 switches from "reading" to "Dyslexia"
$ws1.Range("X2").Value = "Dyslexia"

# Row 3 (the second study entry) is removed entirely
$ws1.Rows("3:3").Delete()

# --- Sheet 2: SNP_entryData ---
$ws2 = $wb.Worksheets.Item("SNP_entryData")

# G2 already holds the text value "1" - use it as a clean source for
# PasteSpecial so the other "1"-as-text cells pick up the same string
# cell (shared string / text type) without touching any formatting/styles.
$textOneCells = @("B2","D2","E2","K2","L2","M2","Q2","R2","S2","V2","W2","Z2")
foreach ($cellRef in $textOneCells) {
    $ws2.Range("G2").Copy()
    $ws2.Range($cellRef).PasteSpecial(-4163)
}

# Plain numeric 1's
$numericOneCells = @("A2","C2","F2","H2","I2","J2","N2","O2","P2","U2","X2","Y2")
foreach ($cellRef in $numericOneCells) {
    $ws2.Range($cellRef).Value = 1
}

# New gene-entry specific text
$ws2.Range("T2").Value = "gene"
$ws2.Range("AA2").Value = "google.com"

# AB2 ("disorder" = "reading") is left untouched.

# Row 3 (the second SNP entry) is removed entirely
$ws2.Rows("3:3").Delete()
